$d = $word.ActiveDocument

# The document ends with a handful of empty paragraphs, the last of which
# carries the "_GoBack" bookmark (left over from the last edit position).
# We need to insert two new paragraphs of text, an empty paragraph, and
# then add text before/after the bookmark inside the final paragraph.

# Insert three new (empty) paragraphs right before the bookmark paragraph,
# re-fetching the "last paragraph" each time so we keep inserting right
# before the real bookmark paragraph rather than a stale reference.
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphBefore()
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphBefore()
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphBefore()

# Re-fetch a fresh reference to the bookmark paragraph (now the last one)
# and walk backwards to the three paragraphs we just inserted.
$bookmarkPara = $d.Paragraphs.Last
$idea = $bookmarkPara.Previous(3)
$bears = $bookmarkPara.Previous(2)
# $blank = $bookmarkPara.Previous(1)   # stays empty

$idea.Range.Text = "Ideia para apresentar:"
$bears.Range.Text = "*Usar ursinhos de pelúcia para demostrar de forma mais simples do método de aprendizagem da IA"

# Add the new text after the bookmark FIRST, while the bookmark paragraph is
# still empty: InsertAfter on a range collapsed at the end of the (empty)
# paragraph lands after the zero-width bookmark marks but still before the
# paragraph mark. Doing this before adding any "before bookmark" text keeps
# the End position unambiguous (once there is real text in the paragraph,
# Range.End refers to the end of that text, which is before the bookmark).
$afterBookmark = $d.Range($bookmarkPara.Range.End, $bookmarkPara.Range.End)
$afterBookmark.InsertAfter("(pouco provável)*")

# Now add the new text before the bookmark, at the (still unchanged) start
# of the paragraph.
$beforeBookmark = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)
$beforeBookmark.InsertBefore("*Usar cabine para fazer teste de turing* ")
